# "single stream data for paper"
#
# The "Calibration from 03/10/2015" block (originally columns L:S on
# Sheet1, with its paired slope/intercept readout in columns R:S) is
# widened to make room for a new single-stream data series: two blank
# columns are inserted immediately to the right of the first block's
# trailing spacer column (J), pushing the second block from L:S to N:U.
#
# This is exactly an Excel "insert 2 columns" at K:L, which also carries
# along the existing J-column spacer formatting for the two new columns,
# shifts every merged range / formula in those columns two letters to the
# right, and grows the sheet dimension from T9 to V9.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert two blank columns at K:L, shifting everything from the old
# column K onward two places to the right (K->M, L->N, ... T->V).
$ws.Columns("K:L").Insert()

# The newly freed spacer cells flanking the (now shifted) data blocks
# are centered, matching the formatting applied next to the other
# calibration blocks.
$ws.Range("J2:K3").HorizontalAlignment = -4108
$ws.Range("J5:K6").HorizontalAlignment = -4108
$ws.Range("J8:K9").HorizontalAlignment = -4108

# Leave the selection where the author ended up after the edit.
[void]$ws.Range("I18").Select()
